$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 38, shifting existing rows 38..67 down to 39..68
$ws.Rows.Item(38).Insert()

# Populate the newly inserted row 38 with the new weekly record
$ws.Cells.Item(38, 1).Value = 9
$ws.Cells.Item(38, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(38, 3).Value = "Metropolitana"
$ws.Cells.Item(38, 4).Value = 44741
$ws.Cells.Item(38, 5).Value = 13
$ws.Cells.Item(38, 6).Value = 100112029
$ws.Cells.Item(38, 7).Value = "Orégano"
$ws.Cells.Item(38, 8).Value = "Sin especificar"
$ws.Cells.Item(38, 9).Value = "Primera"
$ws.Cells.Item(38, 10).Value = 16
$ws.Cells.Item(38, 11).Value = 16000
$ws.Cells.Item(38, 12).Value = 16000
$ws.Cells.Item(38, 13).Value = 16000
$ws.Cells.Item(38, 14).Value = "$/docena de atados"
$ws.Cells.Item(38, 15).Value = "Región Metropolitana"
$ws.Cells.Item(38, 16).Value = 5333
$ws.Cells.Item(38, 17).Value = 3
$ws.Cells.Item(38, 18).Value = "Hortaliza"
